$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Update hyperlink target URLs for rows that survive (F3, F4), and
#        drop the hyperlinks that belonged to rows which are about to be
#        removed (F5:F13). F2 keeps its original target.
$toDeleteLinks = New-Object System.Collections.ArrayList
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$F$3') {
        $hl.Address = "https://www.lancers.jp/work/detail/5415610"
    } elseif ($addr -eq '$F$4') {
        $hl.Address = "https://www.lancers.jp/work/detail/5415615"
    } elseif ($addr -ne '$F$2') {
        [void]$toDeleteLinks.Add($hl)
    }
}
for ($i = $toDeleteLinks.Count - 1; $i -ge 0; $i--) {
    $toDeleteLinks[$i].Delete()
}

# --- 2. Remove rows 5-13 entirely (also shrinks the used range / dimension).
$ws.Range("A5:H13").EntireRow.Delete()

# --- 3. Refresh the "取得日時" timestamp on the remaining rows.
$ws.Range("A2").Value = "2025-10-18 06:24:18"
$ws.Range("A3").Value = "2025-10-18 06:24:18"
$ws.Range("A4").Value = "2025-10-18 06:24:18"

# --- 4. Row 3 now carries the data that used to live in row 8.
$ws.Range("B3").Value = "仮想通貨取引のBOT作成"
$ws.Range("D3").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5415610"
$ws.Range("G3").Value = 118
$ws.Range("H3").Value = "★bot"

# --- 5. Row 4 now carries the data that used to live in row 11 (no H value).
$ws.Range("B4").Value = "【クリエイティブ】Aurora Creative Lab 外注パートナー募集"
$ws.Range("D4").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5415615"
$ws.Range("G4").Value = 18
$ws.Range("H4").ClearContents()

# --- 6. Column B narrows from 49 to 40 characters wide.
#        (The engine's ColumnWidth property is offset by 5/6 of a character
#        from the stored OOXML "width" attribute, so compensate here.)
$ws.Columns.Item(2).ColumnWidth = 40 - 5/6
